# CPT-280 style: customize header row's stylings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Text changes: prefix the locale codes in the header row
#    B1 "zh" -> "Locale: zh", C1 "en" -> "Locale: en"
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Locale: zh"
$ws.Range("C1").Value = "Locale: en"

# ---------------------------------------------------------------------
# 2) Column widths: keep column A at 50, widen B:C to 150
#    (ColumnWidth is in "characters"; the engine adds a fixed 5/6
#    character padding when it serialises to OOXML <col width>, so we
#    compensate by subtracting it up front to land exactly on 150.)
# ---------------------------------------------------------------------
$ws.Columns("B:C").ColumnWidth = 150 - 0.8333333333333334

# ---------------------------------------------------------------------
# 3) Header row height
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 40

# ---------------------------------------------------------------------
# 4) Alignment: center every header cell horizontally (vertical center
#    is already inherited from the existing style).
# ---------------------------------------------------------------------
$ws.Range("A1").HorizontalAlignment = -4108
$hdr = $ws.Range("B1:C1")
$hdr.HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5) Bold + blue font for the B1:C1 "value" header cells
# ---------------------------------------------------------------------
$hdr.Font.Bold = $true
$hdr.Font.Color = 15351613

# ---------------------------------------------------------------------
# 6) Freeze panes at B2 so both the header row and the first (key)
#    column stay pinned.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit applied"
